$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap columns I and K, and columns M and O for data rows 2-25
for ($r = 2; $r -le 25; $r++) {
    $iVal = $ws.Cells.Item($r, 9).Value2   # column I
    $kVal = $ws.Cells.Item($r, 11).Value2  # column K
    $ws.Cells.Item($r, 9).Value = $kVal
    $ws.Cells.Item($r, 11).Value = $iVal

    $mVal = $ws.Cells.Item($r, 13).Value2  # column M
    $oVal = $ws.Cells.Item($r, 15).Value2  # column O
    $ws.Cells.Item($r, 13).Value = $oVal
    $ws.Cells.Item($r, 15).Value = $mVal
}

# Add new columns P and Q with value 2 for rows 2-25
$ws.Range("P2:Q25").Value = 2

# Add new header cells P1 and Q1, matching style of existing headers (O1)
$ws.Range("O1").Copy() | Out-Null
$ws.Range("P1:Q1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15
